$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "updating ssn in persona": replace the placeholder SSNs in column L
# (rows 2-7) with the updated values.
$ws.Range("L2").Value = 238435798
$ws.Range("L3").Value = 238435799
$ws.Range("L4").Value = 238435800
$ws.Range("L5").Value = 238435801
$ws.Range("L6").Value = 238435802
$ws.Range("L7").Value = 238435803

# Leave the selection where the author left it when they saved.
[void]$ws.Range("M9").Select()
